# Update the "Plot" column (A) labels from single-letter codes to full
# descriptive names: C -> Control, F -> Freshwater, S -> Seawater.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    switch ($cell.Value2) {
        "C" { $cell.Value = "Control" }
        "F" { $cell.Value = "Freshwater" }
        "S" { $cell.Value = "Seawater" }
    }
}

# Update the active selection on the sheet to match the edited file.
$ws.Range("A22:A31").Select()
